$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242 (weekly price-log entry), pushing the
# existing rows 242-254 down to 243-255.
$ws.Rows.Item(242).Insert()

$ws.Range("A242").Value = 4
$ws.Range("B242").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C242").Value = "Los Lagos"
$ws.Range("D242").Value = 44610
$ws.Range("E242").Value = 10
$ws.Range("F242").Value = 100114014
$ws.Range("G242").Value = "Betarraga"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 500
$ws.Range("K242").Value = 1000
$ws.Range("L242").Value = 1000
$ws.Range("M242").Value = 1000
$ws.Range("N242").Value = "$/paquete 5 unidades"
$ws.Range("O242").Value = "Región del Maule"
$ws.Range("P242").Value = 200
$ws.Range("Q242").Value = 5
$ws.Range("R242").Value = "Hortaliza"
